# Generate Report for Handoff
# Status moves from "In Translation" to "Ready for handoff" and the
# handoff timestamps advance a few minutes on the Overview sheet and on
# each per-language sheet (zh-cn, de-de). The Status columns also widen
# to fit the new, longer status text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$ws = $wb.Worksheets("Overview")
$ws.Range("E2").Value = "Ready for handoff"   # zh-cn status
$ws.Range("F2").Value = "Ready for handoff"   # de-de status
$ws.Range("G2").Value = "2016-08-26 10:37:40" # Latest HO Xliff Generate Date
$ws.Columns(5).ColumnWidth = 16.25            # widen zh-cn column for new text
$ws.Columns(6).ColumnWidth = 16.25            # widen de-de column for new text

# --- zh-cn sheet --------------------------------------------------------
$ws = $wb.Worksheets("zh-cn")
$ws.Range("C2").Value = "Ready for handoff"   # Status
$ws.Range("H2").Value = "2016-08-26 10:37:36" # Latest Handoff Datetime
$ws.Columns(3).ColumnWidth = 16.25            # widen Status column for new text

# --- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets("de-de")
$ws.Range("C2").Value = "Ready for handoff"   # Status
$ws.Range("H2").Value = "2016-08-26 10:37:40" # Latest Handoff Datetime
$ws.Columns(3).ColumnWidth = 16.25            # widen Status column for new text
